$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title placeholder: collapse the word-by-word runs into a single run.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Characters(1, $titleRange.Length).Text = "Testing custom properties"

# Subtitle placeholder: collapse the runs either side of the two
# existing <a:br/> line breaks, without disturbing the breaks themselves.
$subRange = $s.Shapes.Item(2).TextFrame.TextRange
$subRange.Characters(1, 18).Text = "This is a subtitle"
$subRange.Characters(21, 5).Text = "A. M."
